$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet3")
$ws2 = $wb.Worksheets.Item("Sheet4")

# --- Seed the shared-string table in the exact order the target file uses ---
# (written to an out-of-the-way scratch range, then cleared once every
# target cell below has taken its own reference to the string)
$order = @(
    "Contribuições",
    "Doaçoes",
    "Nome Completo - CPF",
    "Banco/Conta",
    "Receita",
    "Ativo",
    "Despesa",
    "Transporte e Viagens",
    "Propaganda Institucional",
    "Venda de Material para Divulgação",
    "Material para Comercialização",
    "Despesas com alistamento"
)
for ($i = 0; $i -lt $order.Length; $i++) {
    $ws2.Cells.Item(200 + $i, 26).Value = $order[$i]
}

# --- New "ledger example" helper columns (H/I/J) on Sheet4 ---
$ws2.Range("H4").Value = "Receita"
$ws2.Range("I4").Value = "Contribuições"

$ws2.Range("J5").Value = "Nome Completo - CPF"

$ws2.Range("H7").Value = "Receita"
$ws2.Range("I7").Value = "Venda de Material para Divulgação"

$ws2.Range("H9").Value = "Receita"
$ws2.Range("I9").Value = "Doaçoes"

$ws2.Range("J10").Value = "Nome Completo - CPF"

$ws2.Range("H11").Value = "Ativo"
$ws2.Range("I11").Value = "Caixa"

$ws2.Range("J12").Value = "Banco/Conta"

$ws2.Range("H15").Value = "Despesa"
$ws2.Range("I15").Value = "Material de Consumo"

$ws2.Range("H16").Value = "Estoque"
$ws2.Range("I16").Value = "Material para Comercialização"

$ws2.Range("H17").Value = "Despesa"
$ws2.Range("I17").Value = "Aluguéis e Condomínios"

$ws2.Range("H18").Value = "Despesa"
$ws2.Range("I18").Value = "Transporte e Viagens"

$ws2.Range("H19").Value = "Despesa"
$ws2.Range("I19").Value = "Serviços Técnico-Profissionais"

$ws2.Range("H20").Value = "Despesa"
$ws2.Range("I20").Value = "Serviços e Utilidades"

$ws2.Range("H21").Value = "Despesa"
$ws2.Range("I21").Value = "Propaganda Institucional"

$ws2.Range("H22").Value = "Despesa"
$ws2.Range("I22").Value = "Despesas Financeiras"

$ws2.Range("H23").Value = "Despesa"
$ws2.Range("I23").Value = "Despesas com alistamento"

# --- drop the scratch seed cells now that every string has a real owner ---
$ws2.Range("Z200:Z211").ClearContents()

# --- view state: Sheet3 scrolled back to top, selection on D12 ---
$ws1.Activate()
$ws1.Range("A1").Select()
$ws1.Range("D12").Select()

# --- view state: Sheet4 stays the active tab, selection moves to L17 ---
$ws2.Activate()
$ws2.Range("L17").Select()
